# Update countries & provincias Spain
# Applies the data refresh captured in the commit: a handful of country
# rows swap places (because their "Casos totales" crossed over during the
# update) and several rows get refreshed case/recovered/death numbers.
# The "last updated" timestamp in A1 is also bumped by an hour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: refreshed timestamp -------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 31 de Mayo de 2020 a las 11:05"

# --- Rows 24/25: Banglades overtakes Paises Bajos -----------------------
# Row 24 now holds Banglades' (updated) figures, row 25 holds
# Paises Bajos' (updated) figures - the two countries swapped rank.
$ws.Cells.Item(24,1).Value = "Banglades"
$ws.Cells.Item(24,2).Value = 47153
$ws.Cells.Item(24,3).Value = 2545
$ws.Cells.Item(24,4).Value = 9781
$ws.Cells.Item(24,5).Value = 36722
$ws.Cells.Item(24,6).Value = 0
$ws.Cells.Item(24,7).Value = 40
$ws.Cells.Item(24,8).Value = 650

$ws.Cells.Item(25,1).Value = "Paises Bajos"
$ws.Cells.Item(25,2).Value = 46257
$ws.Cells.Item(25,3).Value = 0
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 0
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = 5951

# --- Row 12 (India): refreshed figures ----------------------------------
$ws.Cells.Item(12,2).Value = 182681
$ws.Cells.Item(12,3).Value = 854
$ws.Cells.Item(12,4).Value = 87049
$ws.Cells.Item(12,5).Value = 90446

# --- Row 54 (Barein): refreshed figures ---------------------------------
$ws.Cells.Item(54,5).Value = 4949
$ws.Cells.Item(54,7).Value = 1
$ws.Cells.Item(54,8).Value = 18

# --- Row 62 (Moldavia): refreshed figures -------------------------------
$ws.Cells.Item(62,4).Value = 4581
$ws.Cells.Item(62,5).Value = 3223
$ws.Cells.Item(62,7).Value = 3
$ws.Cells.Item(62,8).Value = 294

# --- Row 82 (Consejo Danes para los Refugiados): refreshed figures -----
$ws.Cells.Item(82,2).Value = 3070
$ws.Cells.Item(82,3).Value = 104
$ws.Cells.Item(82,4).Value = 448
$ws.Cells.Item(82,5).Value = 2551
$ws.Cells.Item(82,7).Value = 2
$ws.Cells.Item(82,8).Value = 71

# --- Row 101 (Sri Lanka): refreshed figures -----------------------------
$ws.Cells.Item(101,4).Value = 801
$ws.Cells.Item(101,5).Value = 809

# --- Row 102 (Eslovaquia): refreshed figures ----------------------------
$ws.Cells.Item(102,4).Value = 1366
$ws.Cells.Item(102,5).Value = 127

# --- Row 112 (Hong Kong): refreshed figures -----------------------------
$ws.Cells.Item(112,2).Value = 1085
$ws.Cells.Item(112,3).Value = 2
$ws.Cells.Item(112,4).Value = 1037
$ws.Cells.Item(112,5).Value = 44

# --- Rows 200/201: Santa Lucia overtakes Belice -------------------------
# Row 200 now holds Santa Lucia's figures, row 201 holds Belice's figures.
$ws.Cells.Item(200,1).Value = "Santa Lucia"
$ws.Cells.Item(200,4).Value = 18
$ws.Cells.Item(200,8).Value = 0

$ws.Cells.Item(201,1).Value = "Belice"
$ws.Cells.Item(201,4).Value = 16
$ws.Cells.Item(201,8).Value = 2

# --- Rows 213/214: Papua Nueva Guinea overtakes Islas Virgenes Britanicas
# Row 213 now holds Papua Nueva Guinea's figures, row 214 holds
# Islas Virgenes Britanicas' figures.
$ws.Cells.Item(213,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213,4).Value = 8
$ws.Cells.Item(213,8).Value = 0

$ws.Cells.Item(214,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214,4).Value = 7
$ws.Cells.Item(214,8).Value = 1
